# Update the "Marking" row (row 11) and "Total" row (row 12) figures
# on the active worksheet to reflect the corrected scoring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right 6 -> 9, Wrong 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 ("Total"): Right 138 -> 207, Wrong -6 -> -4, Max fraction 132/168 -> 203/252
$ws.Range("B12").Value = 207
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "203/252"
